$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new header cells (G1, H1), copying the header formatting
# from the existing F1 header cell (bold, border, centered) ---
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Column H keeps the old "TCV_range" meaning/values that used to live in F.
$ws.Range("H1").Value = "TCV_range"
$ws.Range("H2:H23").Value = "60000-80000"

# Column G is the new "Ny leverandør" (new supplier) column - sparse.
$ws.Range("G1").Value = "Ny leverandør"
$ws.Range("G3").Value = "DataLøn"
$ws.Range("G7").Value = "Zenegy"
$ws.Range("G10").Value = "DataLøn"

# Column F is repurposed from "TCV_range" to "Årsag" (reason for cancellation).
$ws.Range("F1").Value = "Årsag"

$ws.Range("F2").Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Range("F3").Value = "Outsourcing af lønnen (anden leverandør)"
$ws.Range("F4").Value = "Systemet (uddyb i bemærkninger)"
$ws.Range("F5").Value = "Ikke oplyst"
$ws.Range("F6").Value = "Utilfredshed (Ventetid på telefon)"
$ws.Range("F7").Value = "Utilfredshed (Service - uddyb i bemærkninger)"
$ws.Range("F8").Value = "Ikke oplyst"
$ws.Range("F9").Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Range("F10").Value = "Ikke oplyst"
$ws.Range("F11").Value = "Utilfredshed (Service - uddyb i bemærkninger)"
$ws.Range("F12").Value = "Ikke oplyst"
$ws.Range("F13").Value = "Fusionerer med anden virksomhed"
$ws.Range("F14").Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Range("F15").Value = "Ikke oplyst"
$ws.Range("F16").Value = "Ikke oplyst"
$ws.Range("F17").Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Range("F18").Value = "Anden årsag (angiv hvilken i bemærkninger)"
$ws.Range("F19").Value = "Ikke oplyst"
$ws.Range("F20").Value = "Ikke oplyst"
$ws.Range("F21").Value = "Pris"
$ws.Range("F22").Value = "Ikke oplyst"
$ws.Range("F23").Value = "Ikke oplyst"
